# Generate Report for Handback
# Update the "last generated" timestamp columns on each worksheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), row 2
$wsOverview.Range("G2").Value = "2016-08-28 23:04:31"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-08-28 23:04:26"
$wsZhCn.Range("K2").Value = "2016-08-28 23:04:41"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-08-28 23:04:31"
$wsDeDe.Range("K2").Value = "2016-08-28 23:04:48"
